$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Marking" row (row 11): Right marks 4 -> 5, Wrong marks -1 -> -1.2
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = -1.2

# Update "Total" row (row 12): Right total 100 -> 125, Wrong total -1 -> -1.2
$ws.Range("B12").Value = 125
$ws.Range("C12").Value = -1.2

# Update the computed score text in E12
$ws.Range("E12").Value = "123.8/140"
